$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.420.08"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "1.915.11"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.38"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4690"
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2868"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06825"
$ws.Range("E9").Value = "  +3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.63"
$ws.Range("E10").Value = "  +12.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.45"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07726"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "1.884.70"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.275"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6581"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "296.75"
$ws.Range("E16").Value = "  -5.46%  "
$ws.Range("D17").Value = "30.417.82"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007627"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9989"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.92"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").Value = "2.134.79"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9990"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.229"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.204"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.66"
$ws.Range("E25").Value = "  +6.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.304"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.18"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.090"
$ws.Range("E28").Value = "  +6.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1071"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.364"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.162"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.980"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05058"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7369"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02067"
$ws.Range("E36").Value = "  +5.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.744"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.681"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.059"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "109.42"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8714"
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.812"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4254"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "51.46"
$ws.Range("E45").Value = "  +19.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.41"
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.200"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.252"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1212"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.79"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2436"
$ws.Range("E51").Value = "  +10.73%  "
